# Add two newer case records to the top of the "VIC mystery cases" table,
# pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table to its final size (A1:E25) first. Doing this before the
# physical row insert keeps the calculated-column formulas on the rows that
# end up at the bottom (formerly the table's last rows) correctly scoped to
# the table instead of getting rewritten as plain (broken) cell refs.
$lo.Resize($ws.Range("A1:E25"))

# Insert two blank rows at the top of the data (row 2), shifting every
# existing record down by two rows.
$ws.Rows("2:3").Insert()

# Copy the formatting from the row that is now row 4 (the original first
# data row) down into the two new rows so the new records are styled the
# same as the rest of the table.
$ws.Range("A4:E4").Copy()
$ws.Range("A2:E3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the two new records. C3 is written before C2 so the shared
# strings end up in the same order as the source edit (row-3's link first).
$ws.Range("C3").Value = "https://www.dhhs.vic.gov.au/coronavirus-update-for-victoria-29-October-2020"
$ws.Range("C2").Value = "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-3-november-2020"

$ws.Range("A2").Value = 44134
$ws.Range("B2").Value = 3023

$ws.Range("A3").Value = 44130
$ws.Range("B3").Value = 3081

$ws.Range("D2").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
$ws.Range("E2").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"
$ws.Range("D3").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
$ws.Range("E3").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"

# The existing hyperlink lived on the "17-october-2020" record, which was on
# row 5 and is now on row 7 after the insert - move it along with the data.
$ws.Range("C5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C7"), "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020")

# Match the saved selection state of the edited workbook.
$null = $ws.Range("A3").Select()
